# "change server dir names (for deploy and tools)"
#
# The DataNode sheet's server table (A9:H13) gets its rows reordered AND
# the ServerID/Type numbers (columns B and H) updated to new values:
#
#   Name            old ServerID/Type -> new     Port (unchanged, travels with the row)
#   MasterServer_1   3 -> 1     13001
#   WorldServer_1    7 -> 50    17001
#   GameServer_1     6 -> 51    16001
#   LoginServer_1    4 -> 52    14001
#   ProxyServer_1    5 -> 53    15001
#
# New row order (9..13): MasterServer_1, WorldServer_1, GameServer_1,
# LoginServer_1, ProxyServer_1 (per-row formatting - e.g. the thicker
# border on the final data row - must travel with the data, so this is
# done as a physical row reorder rather than independent cell edits).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Stage the five existing rows (9-13) in a scratch area far below the
#    used range so every source row can be read before any destination
#    is overwritten.
$ws.Range("A9:H9").Cut($ws.Range("A101:H101"))    # GameServer_1
$ws.Range("A10:H10").Cut($ws.Range("A102:H102"))  # WorldServer_1
$ws.Range("A11:H11").Cut($ws.Range("A103:H103"))  # ProxyServer_1
$ws.Range("A12:H12").Cut($ws.Range("A104:H104"))  # MasterServer_1
$ws.Range("A13:H13").Cut($ws.Range("A105:H105"))  # LoginServer_1

# 2) Drop them back in the new order, carrying each row's formatting
#    (including the special last-row border) along with its data.
$ws.Range("A104:H104").Cut($ws.Range("A9:H9"))    # -> row 9: MasterServer_1
$ws.Range("A102:H102").Cut($ws.Range("A10:H10"))  # -> row 10: WorldServer_1
$ws.Range("A101:H101").Cut($ws.Range("A11:H11"))  # -> row 11: GameServer_1
$ws.Range("A105:H105").Cut($ws.Range("A12:H12"))  # -> row 12: LoginServer_1
$ws.Range("A103:H103").Cut($ws.Range("A13:H13"))  # -> row 13: ProxyServer_1

# 3) Drop the now-empty scratch area entirely so it doesn't linger as
#    blank formatted rows in the sheet.
$ws.Range("A101:H105").Clear()

# 4) Update the ServerID (B) / Type (H) numbers for each server.
$ws.Range("B9").Value = "1"
$ws.Range("H9").Value = "1"

$ws.Range("B10").Value = "50"
$ws.Range("H10").Value = "50"

$ws.Range("B11").Value = "51"
$ws.Range("H11").Value = "51"

$ws.Range("B12").Value = "52"
$ws.Range("H12").Value = "52"

$ws.Range("B13").Value = "53"
$ws.Range("H13").Value = "53"
